# Rename the "LNIO" condition label to "IONL" in the header row (row 1,
# columns C, E, G, I, K, M, O, Q, S, U) of both worksheets, then restore
# the selections/active sheet exactly as recorded after the edit.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("spike_nb_pyr")
$sheet2 = $wb.Worksheets.Item("spike_nb_int")

foreach ($ws in @($sheet1, $sheet2)) {
    foreach ($cell in @("C1", "E1", "G1", "I1", "K1", "M1", "O1", "Q1", "S1", "U1")) {
        $old = $ws.Range($cell).Value2
        $new = $old -replace "LNIO", "IONL"
        $ws.Range($cell).Value = $new
    }
}

# Update per-sheet selections (set while each sheet is active so the
# selection "sticks" to that sheet).
$sheet1.Activate()
$sheet1.Range("E9").Select() | Out-Null

# Make "spike_nb_int" (second sheet, index 1) the active/selected tab -
# matching the workbook's new activeTab="1".
$sheet2.Activate()
$sheet2.Range("G11").Select() | Out-Null
